# Backfill nine missing trading days (2019-11-18 .. 2019-11-28) that were absent
# between the existing 2019-11-15 and 2019-11-29 rows.
#
# Effect: the sheet grows from A1:I833 to A1:I842 - the old rows 761-833
# (2019-11-29 onward) shift down nine rows to become 770-842, and the newly
# inserted rows 761-769 hold the backfilled OHLCV data below.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteValues = -4163

# Push the existing 2019-11-29..2020-03-13 block (rows 761-833) down nine rows.
$ws.Range("A761:A769").EntireRow.Insert()

# New rows to create, in order, right after row 760 (2019-11-15).
# Columns: row, timestamp(A), date(B), open(E), high(F), low(G), close(H), vol(I)
# id(C)="0186" and name(D)="PTRANS" are constant for every row in this sheet.
$newRows = @(
    @{ Row = 761; Ts = 1574035200; Date = "2019-11-18"; Open = 0.255; High = 0.255; Low = 0.245; Close = 0.255; Vol = 5990000 },
    @{ Row = 762; Ts = 1574121600; Date = "2019-11-19"; Open = 0.255; High = 0.255; Low = 0.245; Close = 0.245; Vol = 2947400 },
    @{ Row = 763; Ts = 1574208000; Date = "2019-11-20"; Open = 0.245; High = 0.25;  Low = 0.245; Close = 0.25;  Vol = 1302600 },
    @{ Row = 764; Ts = 1574294400; Date = "2019-11-21"; Open = 0.25;  High = 0.25;  Low = 0.245; Close = 0.245; Vol = 2218200 },
    @{ Row = 765; Ts = 1574380800; Date = "2019-11-22"; Open = 0.26;  High = 0.265; Low = 0.25;  Close = 0.265; Vol = 49261400 },
    @{ Row = 766; Ts = 1574640000; Date = "2019-11-25"; Open = 0.265; High = 0.27;  Low = 0.26;  Close = 0.265; Vol = 9489200 },
    @{ Row = 767; Ts = 1574726400; Date = "2019-11-26"; Open = 0.265; High = 0.27;  Low = 0.26;  Close = 0.265; Vol = 7531400 },
    @{ Row = 768; Ts = 1574812800; Date = "2019-11-27"; Open = 0.265; High = 0.27;  Low = 0.26;  Close = 0.26;  Vol = 5511300 },
    @{ Row = 769; Ts = 1574899200; Date = "2019-11-28"; Open = 0.26;  High = 0.265; Low = 0.26;  Close = 0.265; Vol = 2368500 }
)

foreach ($nr in $newRows) {
    $r = $nr.Row

    # Clone row 760 (A:I) into the new row first. This brings over the id/name
    # text columns (C, D) - identical on every row of this sheet - as plain
    # text cells with no extra number formatting, matching the rest of the
    # column rather than introducing a one-off style.
    $ws.Range("A760:I760").Copy($ws.Range("A" + $r + ":I" + $r))

    $ws.Range("A" + $r).Value2 = $nr.Ts

    # Write the date as a text formula, then convert it to a plain literal by
    # copying the cell onto itself with values-only paste. This keeps the
    # date stored as text (like "2019-11-18"), not as an auto-converted date
    # serial number, without leaving a formula behind or needing a quote-
    # prefixed/text-formatted cell style.
    $ws.Range("B" + $r).Formula = '="' + $nr.Date + '"'
    $ws.Range("B" + $r).Copy()
    $ws.Range("B" + $r).PasteSpecial($xlPasteValues)

    $ws.Range("E" + $r).Value2 = $nr.Open
    $ws.Range("F" + $r).Value2 = $nr.High
    $ws.Range("G" + $r).Value2 = $nr.Low
    $ws.Range("H" + $r).Value2 = $nr.Close
    $ws.Range("I" + $r).Value2 = $nr.Vol
}
